$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 189, pushing the existing
# rows 189-230 down to 191-232 (unchanged), and populate the two newly
# created rows (189-190) with their own data.
$ws.Rows("189:190").Insert()

# New row 189
$ws.Range("A189").Value = 7
$ws.Range("B189").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C189").Value = "Ñuble"
$ws.Range("D189").Value = 44694
$ws.Range("E189").Value = 16
$ws.Range("F189").Value = 100112043
$ws.Range("G189").Value = "Pepino ensalada"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 60
$ws.Range("K189").Value = 14500
$ws.Range("L189").Value = 15000
$ws.Range("M189").Value = 14750
$ws.Range("N189").Value = "$/caja 60 unidades"
$ws.Range("O189").Value = "Región de Arica y Parinacota"
$ws.Range("P189").Value = 246
$ws.Range("Q189").Value = 60
$ws.Range("R189").Value = "Hortaliza"

# New row 190
$ws.Range("A190").Value = 7
$ws.Range("B190").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C190").Value = "Ñuble"
$ws.Range("D190").Value = 44694
$ws.Range("E190").Value = 16
$ws.Range("F190").Value = 100112043
$ws.Range("G190").Value = "Pepino ensalada"
$ws.Range("H190").Value = "Sin especificar"
$ws.Range("I190").Value = "Primera"
$ws.Range("J190").Value = 120
$ws.Range("K190").Value = 15000
$ws.Range("L190").Value = 16000
$ws.Range("M190").Value = 15500
$ws.Range("N190").Value = "$/caja 80 unidades"
$ws.Range("O190").Value = "Región del Maule"
$ws.Range("P190").Value = 194
$ws.Range("Q190").Value = 80
$ws.Range("R190").Value = "Hortaliza"
